# --- Rename the worksheet tab: "CoA Import" -> "Chart of Accounts" ---
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Chart of Accounts"

# --- Header row (row 1): re-case the column headers ---
$ws.Cells.Item(1,1).Value = "Code"
$ws.Cells.Item(1,2).Value = "Name"
$ws.Cells.Item(1,3).Value = "Type"
$ws.Cells.Item(1,4).Value = "Subtype"
$ws.Cells.Item(1,5).Value = "IsSubledger"
$ws.Cells.Item(1,6).Value = "SubledgerType"
$ws.Cells.Item(1,7).Value = "Active"
$ws.Cells.Item(1,8).Value = "Description"
$ws.Cells.Item(1,9).Value = "ParentCode"

# --- Replace the sample data (rows 2-11) with the new Chart-of-Accounts rows ---
# Columns: Code, Name, Type, Subtype, IsSubledger, SubledgerType, Active, Description, ParentCode
$data = @(
  @("1000","Test Cash","ASSET","Current Asset","false","","true","Cash on hand",""),
  @("1100","Test Bank Account","ASSET","Current Asset","false","","true","Primary bank account","1000"),
  @("1200","Test Accounts Receivable","ASSET","Current Asset","true","Customer","true","Amounts owed by customers",""),
  @("2000","Test Accounts Payable","LIABILITY","Current Liability","true","Vendor","true","Amounts owed to vendors",""),
  @("2100","Test Credit Card","LIABILITY","Current Liability","false","","true","Business credit card","2000"),
  @("3000","Test Owner Equity","EQUITY","Equity","false","","true","Owner investment",""),
  @("4000","Test Sales Revenue","REVENUE","Operating Revenue","false","","true","Revenue from sales",""),
  @("4100","Test Service Revenue","REVENUE","Operating Revenue","false","","true","Revenue from services","4000"),
  @("5000","Test Rent Expense","EXPENSE","Operating Expense","false","","true","Office rent",""),
  @("5100","Test Utilities Expense","EXPENSE","Operating Expense","false","","true","Electricity, water, etc.","5000")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 1
        $v = $values[$j]
        $isNumeric = $v -match '^[0-9]+$'
        $isBoolLike = ($v -eq "true") -or ($v -eq "false")
        $isBlank = ($v -eq "")
        if ($isNumeric -or $isBoolLike -or $isBlank) {
            # These source values are stored as literal TEXT in the workbook
            # (not real numbers/booleans/empty cells), so numeric-looking
            # codes ("1000"), blank cells ("") and the "true"/"false" flags
            # need a leading apostrophe to stop Excel auto-converting them
            # into a number / boolean / empty cell.
            $ws.Cells.Item($row, $col).Value = "'" + $v
        } else {
            $ws.Cells.Item($row, $col).Value = $v
        }
    }
}
